$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix accented / capitalisation typos in existing island names ---
$ws.Range("A2").Value = "L’île de la vieille dame"
$ws.Range("A5").Value = "L’île du Piton"
$ws.Range("A6").Value = "L’île des francs"

# --- Add the new "situation" column (I) ---
$ws.Range("I1").Value = "situation"

$sud = "Sud"
$nord = "Nord"

$ws.Range("I2").Value  = $sud
$ws.Range("I3").Value  = $sud
$ws.Range("I4").Value  = $sud
$ws.Range("I5").Value  = $sud
$ws.Range("I6").Value  = $nord
$ws.Range("I7").Value  = $nord
$ws.Range("I8").Value  = $sud
$ws.Range("I9").Value  = $nord
$ws.Range("I10").Value = $nord
$ws.Range("I11").Value = $nord
$ws.Range("I12").Value = $nord
$ws.Range("I13").Value = $nord
$ws.Range("I14").Value = $sud
$ws.Range("I15").Value = $sud
$ws.Range("I16").Value = $sud
$ws.Range("I17").Value = $sud
$ws.Range("I18").Value = $sud

# --- Match the selection left behind by the author (last touched cell) ---
$ws.Range("I18").Select()
